$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# The handback transform failed for the 90fa7d6b... file, so its status
# changes from "Ready for handoff" to "Handback transform failed" across
# every sheet that surfaces that status (Overview rollup + per-locale rows).
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Error Detail column (P) needs to be wide enough to show the new message.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.14
$wsDeDe.Columns.Item(16).ColumnWidth = 39.14

# Populate the new Error Detail text explaining the handback/handoff file
# name mismatch, per locale.
$wsZhCn.Range("P3").Value = "Handback file name: xzmmmwfv.400 is different with handoff file name: 90fa7d6b-8127-4847-8776-ceae8553ab55.80d2f503eeef83f237452c36df412b088ebc8da6.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: xzmmmwfv.400 is different with handoff file name: 90fa7d6b-8127-4847-8776-ceae8553ab55.80d2f503eeef83f237452c36df412b088ebc8da6.de-de."
